# Edit script: split "Миколаївська область" -> "Миколаївська" and add
# responsible person / phone columns (split out of existing column I),
# pushing the "Присутній в меддаті?" column from J to L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet and fix the region name cells (B2:B24)
$ws.Name = "Миколаївська"

for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 2).Value = "Миколаївська"
}

# 2. Insert two new columns at J (pushes the existing J column - "Присутній
#    в меддаті?" - to L) so we can host the new "ПІБ відповідальної особи"
#    and "Телефон" columns.
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(10).Insert()

# 3. New column headers
$ws.Cells.Item(1, 10).Value = "ПІБ відповідальної особи"
$ws.Cells.Item(1, 11).Value = "Телефон"

# Phone numbers must stay text (leading zeros must be preserved).
$ws.Range("K2:K24").NumberFormat = "@"

# 4. Fill in the responsible person's name (J) and phone number (K) per row,
#    derived from the combined "contact person" column (I).
$contacts = @{
    2  = @("Петриченко С.Р.", "0509501386")
    3  = @("Дворецька Н.А.", "06684599911")
    4  = @("Кожухар О.М.", "0500702401")
    5  = @("Домущей Н.А.", "0634649363")
    6  = @("Блоховита О.В.", "0666474378")
    7  = @("Філіпова Л.Ф.", "097 562 45 01")
    8  = @("Мальцева О.М.", "0508472002")
    9  = @("Загляда Л.В.", "0679050435")
    10 = @("Калинка Г.В.", "0988468827")
    11 = @("Омеленчук Г.А.", "0689975280")
    12 = @("Петриченко С.Р.", "0509501386")
    13 = @("Косенчук О. М.", "0682379995")
    14 = @("Левченко Л.В.", "0675926126")
    15 = @("Слюсаренко О.В.", "0960773650")
    16 = @("Черненко О.Б.", "0962044964")
    17 = @("Клиш Л.О.", "0995019056")
    18 = @("Волошина Л.О.", "0677416909")
    19 = @("Стульник І. О", "0982835217")
    20 = @("Гаврилюк Д.О.", "0985828514")
    21 = @("Лєнчевська Л.В.", "0957306981")
    22 = @("Левченко Л.В.", "0675926126")
    23 = @("Трушковська Т. С.", "0982196763")
    24 = @("Шіллер Л.О.", "0960242213")
}

foreach ($r in $contacts.Keys) {
    $pair = $contacts[$r]
    $ws.Cells.Item($r, 10).Value = $pair[0]
    $ws.Cells.Item($r, 11).Value = $pair[1]
}
